# Generate Report for Handoff
# Replaces the two "handed back" sample rows with freshly handed-off files:
# the source documents, hashes, dates and status all move forward, and the
# now-irrelevant "Latest Target File" / "Latest Handback File" columns (and
# their hyperlinks) are cleared because these files have not been handed
# back yet.

$wb = $excel.ActiveWorkbook

# ---- new identifiers -------------------------------------------------
$uuid1  = "0646e59e-8ec3-4821-a9de-21864af14c32"
$uuid2  = "ffff17bc4d83-0e98-4527-855f-5854249b09e1"
$hash   = "d38632e3469738437b3153b9189d010b22a7957f"

$status        = "Ready for handoff"
$overviewDate  = "2016-49-20 00:49:17"
$zhcnDate      = "2016-03-20 00:49:14"
$dedeDate      = "2016-03-20 00:49:17"
$emptyDate     = "0001-01-01 00:00:00"

$md1 = "$uuid1.md"
$md2 = "$uuid2.md"
$zhcnXlf = "$uuid1.$hash.zh-cn.xlf"
$dedeXlf = "$uuid1.$hash.de-de.xlf"

# ---- hyperlink target URLs (kept consistent with the existing scheme) -
$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/0b7a7e77138c1f6f859b0c9a076eff12bef6ecf3/e2e/$md1"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/0b7a7e77138c1f6f859b0c9a076eff12bef6ecf3/e2e/$md2"
$zhcnHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2b44fb760139019c3e7b68353811d834b3af4dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhcnXlf"
$dedeHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e10cc10f777b075f208c7e22e5aa2f5fdeaaf812/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$dedeXlf"

# =========================================================================
# Overview sheet
# =========================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $md1
$ws.Range("B2").Value = $status
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $overviewDate

$ws.Range("A3").Value = $md2
$ws.Range("B3").Value = $status
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $overviewDate

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $md2)

# =========================================================================
# zh-cn sheet
# =========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $md1
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $zhcnXlf
$ws.Range("E2").Value = $zhcnDate
$ws.Range("F2").Clear()
$ws.Range("G2").Clear()
$ws.Range("H2").Value = $emptyDate
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = $md2
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $zhcnXlf
$ws.Range("E3").Value = $zhcnDate
$ws.Range("F3").Clear()
$ws.Range("G3").Clear()
$ws.Range("H3").Value = $emptyDate
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl1, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $zhcnHandoffUrl, "", "", $zhcnXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $md2)
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl2, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $zhcnHandoffUrl, "", "", $zhcnXlf)

# =========================================================================
# de-de sheet
# =========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $md1
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $dedeXlf
$ws.Range("E2").Value = $dedeDate
$ws.Range("F2").Clear()
$ws.Range("G2").Clear()
$ws.Range("H2").Value = $emptyDate
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = $md2
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $dedeXlf
$ws.Range("E3").Value = $dedeDate
$ws.Range("F3").Clear()
$ws.Range("G3").Clear()
$ws.Range("H3").Value = $emptyDate
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl1, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $dedeHandoffUrl, "", "", $dedeXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $md2)
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl2, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $dedeHandoffUrl, "", "", $dedeXlf)
